$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 9348.75
$ws.Range("J88").Value = 9831.666999999999
$ws.Range("L88").Value = 9831.666999999999
$ws.Range("N88").Value = -10643.667

$ws.Range("H91").Value = 9348.75
$ws.Range("J91").Value = 9831.666999999999
$ws.Range("L91").Value = 9831.666999999999
$ws.Range("N91").Value = -12639.667

$ws.Range("H112").Value = 2668.5557
$ws.Range("J112").Value = 2707.8823
$ws.Range("L112").Value = 8123.646900000001
$ws.Range("N112").Value = -10339.6469

$ws.Range("H116").Value = 32171
$ws.Range("I116").Value = 10951.667
$ws.Range("K116").Value = 10951.667
$ws.Range("M116").Value = -7509.666999999999

$ws.Range("H129").Value = 3824.75
$ws.Range("I129").Value = 5599.5
$ws.Range("J129").Value = 2050
$ws.Range("K129").Value = 16798.5
$ws.Range("L129").Value = 6150
$ws.Range("M129").Value = -11798.5
$ws.Range("N129").Value = -16150

$ws.Range("H137").Value = 2516.75
$ws.Range("I137").Value = 2722.3333
$ws.Range("K137").Value = 8166.999899999999
$ws.Range("M137").Value = -5616.999899999999

$ws.Range("H138").Value = 4217.5713
$ws.Range("I138").Value = 2274.25
$ws.Range("K138").Value = 6822.75
$ws.Range("M138").Value = -1682.75

$ws.Range("H141").Value = 6000
$ws.Range("I141").Value = 6000
$ws.Range("K141").Value = 18000
$ws.Range("M141").Value = -12820

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6012
$ws.Range("I61").Value = 6012
$ws.Range("K61").Value = 6012
$ws.Range("M61").Value = -5800

$ws.Range("H110").Value = 8599.75
$ws.Range("J110").Value = 13199.5
$ws.Range("L110").Value = 13199.5
$ws.Range("N110").Value = -17289.5

$ws.Range("H132").Value = 3834.7273
$ws.Range("I132").Value = 3880.6667
$ws.Range("K132").Value = 11642.0001
$ws.Range("M132").Value = -9112.000100000001

$ws.Range("H136").Value = 6012
$ws.Range("I136").Value = 6012
$ws.Range("K136").Value = 18036
$ws.Range("M136").Value = -15486

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1495
$ws.Range("I107").Value = 1495
$ws.Range("K107").Value = 1495
$ws.Range("M107").Value = 425

$ws.Range("H134").Value = 10521.714
$ws.Range("I134").Value = 7844.25
$ws.Range("K134").Value = 23532.75
$ws.Range("M134").Value = -20997.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 378.89474
$ws.Range("J22").Value = 299.5
$ws.Range("L22").Value = 299.5
$ws.Range("N22").Value = -999.5

$ws.Range("H58").Value = 3537.3635
$ws.Range("I58").Value = 3527.7368
$ws.Range("J58").Value = 3598.3333
$ws.Range("K58").Value = 3527.7368
$ws.Range("L58").Value = 3598.3333
$ws.Range("M58").Value = -3324.7368
$ws.Range("N58").Value = -4004.3333

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").Value = $null

$ws.Range("H132").Value = 5017.125
$ws.Range("I132").Value = 4784.75
$ws.Range("J132").Value = 5249.5
$ws.Range("K132").Value = 14354.25
$ws.Range("L132").Value = 15748.5
$ws.Range("M132").Value = -11824.25
$ws.Range("N132").Value = -20808.5

$ws.Range("H134").Value = 7237.3335
$ws.Range("I134").Value = 8244.799999999999
$ws.Range("K134").Value = 24734.4
$ws.Range("M134").Value = -22199.4

$ws.Range("H136").Value = 3537.3635
$ws.Range("I136").Value = 3527.7368
$ws.Range("J136").Value = 3598.3333
$ws.Range("K136").Value = 10583.2104
$ws.Range("L136").Value = 10794.9999
$ws.Range("M136").Value = -8033.2104
$ws.Range("N136").Value = -15894.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 8498.75
$ws.Range("J62").Value = 8498.75
$ws.Range("L62").Value = 25496.25
$ws.Range("N62").Value = -26868.25

$ws.Range("H65").Value = 8498.75
$ws.Range("J65").Value = 8498.75
$ws.Range("L65").Value = 76488.75
$ws.Range("N65").Value = -83352.75

$ws.Range("H80").Value = 11999.75
$ws.Range("J80").Value = 11333
$ws.Range("L80").Value = 33999
$ws.Range("N80").Value = -35871

$ws.Range("H83").Value = 11999.75
$ws.Range("J83").Value = 11333
$ws.Range("L83").Value = 101997
$ws.Range("N83").Value = -111357

$ws.Range("H100").Value = 6333.3335
$ws.Range("J100").Value = 4500
$ws.Range("L100").Value = 13500
$ws.Range("N100").Value = -15122

$ws.Range("H130").Value = 1000
$ws.Range("I130").Value = 1000
$ws.Range("K130").Value = 3000
$ws.Range("M130").Value = 2020

$ws.Range("H132").Value = 3496.6667
$ws.Range("J132").Value = 2993.3333
$ws.Range("L132").Value = 26939.9997
$ws.Range("N132").Value = -31999.9997

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").Value = $null

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").Value = $null

$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").Value = $null

$ws.Range("H132").Value = 3898.75
$ws.Range("I132").Value = 3173
$ws.Range("K132").Value = 9519
$ws.Range("M132").Value = -6989

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 12591.923
$ws.Range("J22").Value = 9857
$ws.Range("L22").Value = 9857
$ws.Range("N22").Value = -10447

$ws.Range("H27").Value = 12591.923
$ws.Range("J27").Value = 9857
$ws.Range("L27").Value = 9857
$ws.Range("N27").Value = -10071

$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = $null
$ws.Range("N132").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 5000000
$ws.Range("I3").Value = 5000000
$ws.Range("M3").Value = -4999886

$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = $null
$ws.Range("N6").Value = $null

$ws.Range("H74").Value = 42488.5
$ws.Range("J74").Value = 39978
$ws.Range("L74").Value = 39978
$ws.Range("N74").Value = -41850

$ws.Range("H77").Value = 42488.5
$ws.Range("J77").Value = 39978
$ws.Range("L77").Value = 119934
$ws.Range("N77").Value = -129294

$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = $null
$ws.Range("N132").Value = $null

$ws.Range("H136").Value = 4500.375
$ws.Range("I136").Value = 4500.375
$ws.Range("K136").Value = 13501.125
$ws.Range("M136").Value = -10951.125
